$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R: "Quantity" header (matches style of existing black-header
# cells like K1/L1 but with a pure-white, non-tinted font) and a numeric
# value of 4 underneath (matches the numeric-cell style used by N2).

# Header cell R1 - copy the format of an existing dark header cell (K1)
# so it reuses the same fill (black) + base font, then flip the font to
# a pure white theme color (no tint) to get the distinct style Excel
# generated for this new header.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("R1").PasteSpecial(-4122) | Out-Null
$ws.Range("R1").Value = "Quantity"
$ws.Range("R1").Font.ThemeColor = 2

# Data cell R2 - copy the format of N2 (center-aligned numeric style)
# and set the quantity value.
$ws.Range("N2").Copy() | Out-Null
$ws.Range("R2").PasteSpecial(-4122) | Out-Null
$ws.Range("R2").Value = 4

$excel.CutCopyMode = 0

# Update the view: scroll window / selection to match the saved state.
$ws.Range("R7").Select() | Out-Null
